$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2 header cell: alignment becomes left / top ----------------------------
$ws.Range("A2").HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Range("A2").VerticalAlignment = -4160     # xlVAlignTop

# --- Row 21: "Nouredine Messalti" entry becomes the bare "Amir" entry --------
$ws.Range("A21").Value = "Amir"
$ws.Range("B21").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

# --- Rows 22-26: brand-new player-name-only rows -----------------------------
$ws.Range("A22").Value = "Naim Dhib"
$ws.Range("A22").HorizontalAlignment = -4108  # xlHAlignCenter (no vertical align)

$ws.Range("A23").Value = "Wael Fareh"
$ws.Range("A23").HorizontalAlignment = -4108  # xlHAlignCenter
$ws.Range("A23").VerticalAlignment = -4108    # xlVAlignCenter

$ws.Range("A24").Value = "Yoan Zouma"
$ws.Range("A24").HorizontalAlignment = -4108
$ws.Range("A24").VerticalAlignment = -4108

$ws.Range("A25").Value = "Ilyes Bougahnmi"
$ws.Range("A25").HorizontalAlignment = -4108
$ws.Range("A25").VerticalAlignment = -4108

$ws.Range("A26").Value = "Oumar"
$ws.Range("A26").HorizontalAlignment = -4108
$ws.Range("A26").VerticalAlignment = -4108

# --- Conditional formatting on A21:A26 (status legend colouring) ------------
# Excel inserts each new rule at the top of the stack (SetFirstPriority),
# so we create them in the reverse of the final priority order.
$rng = $ws.Range("A21:A26")

$cf = $rng.FormatConditions.Add(1, 3, '"NN"')
$cf.Font.Color = 16777215
$cf.Interior.Color = 0
$cf.SetFirstPriority()

$cf = $rng.FormatConditions.Add(1, 3, '"NN"')
$cf.Font.Color = 16777215
$cf.Interior.Color = 4272139
$cf.SetFirstPriority()

$cf = $rng.FormatConditions.Add(1, 3, '"OK"')
$cf.Font.Color = 16777215
$cf.Interior.Color = 5287936
$cf.SetFirstPriority()

$cf = $rng.FormatConditions.Add(1, 3, '"RENFO/TEK"')
$cf.Font.Color = 16777215
$cf.Interior.Color = 15773696
$cf.SetFirstPriority()

$cf = $rng.FormatConditions.Add(1, 3, '"P"')
$cf.Font.Color = 16777215
$cf.Interior.Color = 49407
$cf.SetFirstPriority()

$cf = $rng.FormatConditions.Add(1, 3, '"B"')
$cf.Font.Color = 16777215
$cf.Interior.Color = 255
$cf.SetFirstPriority()

$cf = $rng.FormatConditions.Add(1, 3, '"R"')
$cf.Interior.Color = 16777215
$cf.SetFirstPriority()

# --- Selection cosmetics ------------------------------------------------------
$ws.Range("E23").Select()
